# Generate Report for Handback
# The handback for 90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md has completed,
# so the report is regenerated: that file's row moves to the top (row 2)
# with its status/handback info updated, and ba774427-... moves to row 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 -> 90d90118 file, Row 3 -> ba774427 file
$ov.Cells.Item(2,1).Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$ov.Cells.Item(2,3).Value = ".md"
$ov.Cells.Item(2,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(2,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(2,7).Value = "2016-09-07 08:09:54"

$ov.Cells.Item(3,1).Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$ov.Cells.Item(3,3).Value = ".md"
$ov.Cells.Item(3,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(3,7).Value = "2016-09-07 08:07:50"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "e2e\90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2414d51438a2cb3e38ec85a5b8e4477dc628d5d/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "e2e\ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Cells.Item(2,1).Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$zh.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(2,7).Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.zh-cn.xlf"
$zh.Cells.Item(2,8).Value = "2016-09-07 08:09:42"
$zh.Cells.Item(2,10).Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.zh-cn.xlf"
$zh.Cells.Item(2,11).Value = "2016-09-07 08:10:45"
$zh.Cells.Item(2,13).Value = "True"
$zh.Cells.Item(2,16).Value = ""

$zh.Cells.Item(3,1).Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$zh.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(3,7).Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.zh-cn.xlf"
$zh.Cells.Item(3,8).Value = "2016-09-07 08:07:33"
$zh.Cells.Item(3,10).Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.zh-cn.xlf"
$zh.Cells.Item(3,11).Value = "2016-09-07 08:08:41"
$zh.Cells.Item(3,13).Value = "True"
$zh.Cells.Item(3,16).Value = ""

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/314214a3afac75c4f127c34855b2340849a9bf7c/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2414d51438a2cb3e38ec85a5b8e4477dc628d5d/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/314214a3afac75c4f127c34855b2340849a9bf7c/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null

$zh.Columns.Item(16).ColumnWidth = 12.8

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Cells.Item(2,1).Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$de.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$de.Cells.Item(2,7).Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.de-de.xlf"
$de.Cells.Item(2,8).Value = "2016-09-07 08:09:54"
$de.Cells.Item(2,10).Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.de-de.xlf"
$de.Cells.Item(2,11).Value = "2016-09-07 08:11:09"
$de.Cells.Item(2,13).Value = "True"
$de.Cells.Item(2,16).Value = ""

$de.Cells.Item(3,1).Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$de.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$de.Cells.Item(3,7).Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.de-de.xlf"
$de.Cells.Item(3,8).Value = "2016-09-07 08:07:50"
$de.Cells.Item(3,10).Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.de-de.xlf"
$de.Cells.Item(3,11).Value = "2016-09-07 08:08:59"
$de.Cells.Item(3,13).Value = "True"
$de.Cells.Item(3,16).Value = ""

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/83677caf5d7e6a5e691a9702d99023db283b9a9b/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2414d51438a2cb3e38ec85a5b8e4477dc628d5d/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/83677caf5d7e6a5e691a9702d99023db283b9a9b/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null

$de.Columns.Item(16).ColumnWidth = 12.8
